# Importar RRegistros por nombre de columna
#
# The header row labels are renamed from human-friendly / accented Spanish
# text to plain "code-friendly" identifiers (no spaces/accents) so the
# importer can bind columns by name. Row/column data below the header is
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "NombreSociedad"
$ws.Range("C1").Value = "Ceco"
$ws.Range("Q1").Value = "Area"
$ws.Range("S1").Value = "Guia"
$ws.Range("T1").Value = "FechaSalida"
$ws.Range("U1").Value = "Antiguedad"
$ws.Range("Z1").Value = "FechaActualiza"

# View tweaks that came along with the same save: scroll the frozen-free
# sheet one column to the right and move the active selection to the top
# of the newly relevant column.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$ws.Range("W1").Select()

# New column (Z) got an explicit width in the same edit.
$ws.Columns.Item(26).ColumnWidth = 17.67
